# Commit: table style swap across the deck's three tables.
#
# The source table style GUID ({ED979CD3-0F55-400F-9137-DD503CFAD07B}) is the
# presentation's locally-defined "Table_0" style (see ppt/tableStyles.xml).
# It is replaced everywhere it is used by the built-in PowerPoint table style
# {816C63EA-82A4-42CE-8204-D9BE9A587E19}.

$oldStyleId = "{ED979CD3-0F55-400F-9137-DD503CFAD07B}"
$newStyleId = "{816C63EA-82A4-42CE-8204-D9BE9A587E19}"

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}
